$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 7369.125
$ws.Range("I41").Value = 723.8182
$ws.Range("J41").Value = 21988.8
$ws.Range("K41").Value = 723.8182
$ws.Range("L41").Value = 21988.8
$ws.Range("M41").Value = -283.8182
$ws.Range("N41").Value = -22868.8

$ws.Range("H62").Value = 4442.857
$ws.Range("I62").Value = 3275
$ws.Range("K62").Value = 3275
$ws.Range("M62").Value = -2651

$ws.Range("H65").Value = 4442.857
$ws.Range("I65").Value = 3275
$ws.Range("K65").Value = 16375
$ws.Range("M65").Value = -13255

$ws.Range("H98").Value = 820.5
$ws.Range("I98").Value = 799.8889
$ws.Range("K98").Value = 799.8889
$ws.Range("M98").Value = 698.1111

$ws.Range("H116").Value = 3239.8
$ws.Range("I116").Value = 2947.8235
$ws.Range("K116").Value = 2947.8235
$ws.Range("M116").Value = 494.1765

$ws.Range("H122").Value = 820.5
$ws.Range("I122").Value = 799.8889
$ws.Range("K122").Value = 2399.6667
$ws.Range("M122").Value = 50.33329999999978

$ws.Range("H125").Value = 1306.9
$ws.Range("I125").Value = 1156.8334
$ws.Range("K125").Value = 10411.5006
$ws.Range("M125").Value = -7951.500599999999

$ws.Range("H132").Value = 14070.743
$ws.Range("I132").Value = 1281.0938
$ws.Range("J132").Value = 72537.71000000001
$ws.Range("K132").Value = 3843.2814
$ws.Range("L132").Value = 217613.13
$ws.Range("M132").Value = -1313.2814
$ws.Range("N132").Value = -222673.13

$ws.Range("H137").Value = 2977.5667
$ws.Range("I137").Value = 2778.9048
$ws.Range("K137").Value = 8336.714399999999
$ws.Range("M137").Value = -5786.714399999999

$ws.Range("H138").Value = 2764.05
$ws.Range("I138").Value = 1344.56
$ws.Range("J138").Value = 3409.2727
$ws.Range("K138").Value = 4033.68
$ws.Range("L138").Value = 10227.8181
$ws.Range("M138").Value = 1106.32
$ws.Range("N138").Value = -20507.8181

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 10477.613
$ws.Range("I2").Value = 14773.238
$ws.Range("J2").Value = 1456.8
$ws.Range("K2").Value = 14773.238
$ws.Range("L2").Value = 1456.8
$ws.Range("M2").Value = -14660.238
$ws.Range("N2").Value = -1682.8

$ws.Range("H32").Value = 9539.583000000001
$ws.Range("I32").Value = 8833.391
$ws.Range("K32").Value = 8833.391
$ws.Range("M32").Value = -8546.391

$ws.Range("H63").Value = 3408.0833
$ws.Range("J63").Value = 7000
$ws.Range("L63").Value = 7000
$ws.Range("N63").Value = -8372

$ws.Range("H66").Value = 3408.0833
$ws.Range("J66").Value = 7000
$ws.Range("L66").Value = 35000
$ws.Range("N66").Value = -41864

$ws.Range("H74").Value = 2017.1666
$ws.Range("I74").Value = 2130.1428
$ws.Range("J74").Value = 1621.75
$ws.Range("K74").Value = 2130.1428
$ws.Range("L74").Value = 1621.75
$ws.Range("M74").Value = -1256.1428
$ws.Range("N74").Value = -3369.75

$ws.Range("H77").Value = 2017.1666
$ws.Range("I77").Value = 2130.1428
$ws.Range("J77").Value = 1621.75
$ws.Range("K77").Value = 10650.714
$ws.Range("L77").Value = 8108.75
$ws.Range("M77").Value = -6282.714
$ws.Range("N77").Value = -16844.75

$ws.Range("H116").Value = 10477.613
$ws.Range("I116").Value = 14773.238
$ws.Range("J116").Value = 1456.8
$ws.Range("K116").Value = 14773.238
$ws.Range("L116").Value = 1456.8
$ws.Range("M116").Value = -12479.238
$ws.Range("N116").Value = -6044.8

$ws.Range("H122").Value = 3991
$ws.Range("I122").Value = 1955.1852
$ws.Range("J122").Value = 9487.700000000001
$ws.Range("K122").Value = 5865.5556
$ws.Range("L122").Value = 28463.1
$ws.Range("M122").Value = -3415.5556
$ws.Range("N122").Value = -33363.10000000001

$ws.Range("H132").Value = 3474.3333
$ws.Range("I132").Value = 3223.0625
$ws.Range("K132").Value = 9669.1875
$ws.Range("M132").Value = -7139.1875

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 10477.613
$ws.Range("I3").Value = 14773.238
$ws.Range("J3").Value = 1456.8
$ws.Range("K3").Value = 14773.238
$ws.Range("L3").Value = 1456.8
$ws.Range("M3").Value = -14659.238
$ws.Range("N3").Value = -1684.8

$ws.Range("H20").Value = 1501.4
$ws.Range("I20").Value = 1505.1
$ws.Range("K20").Value = 1505.1
$ws.Range("M20").Value = -1258.1

$ws.Range("H125").Value = 49200
$ws.Range("J125").Value = 49200
$ws.Range("L125").Value = 49200
$ws.Range("N125").Value = -59040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3090.4412
$ws.Range("I31").Value = 1333.0952
$ws.Range("K31").Value = 1333.0952
$ws.Range("M31").Value = -1038.0952

$ws.Range("H34").Value = 3090.4412
$ws.Range("I34").Value = 1333.0952
$ws.Range("K34").Value = 1333.0952
$ws.Range("M34").Value = -1131.0952

$ws.Range("H99").Value = 5349440.5
$ws.Range("J99").Value = 5562193
$ws.Range("L99").Value = 5562193
$ws.Range("N99").Value = -5565189

$ws.Range("H126").Value = 5349440.5
$ws.Range("J126").Value = 5562193
$ws.Range("L126").Value = 16686579
$ws.Range("N126").Value = -16691519

$ws.Range("H132").Value = 3886.1924
$ws.Range("I132").Value = 2909.682
$ws.Range("J132").Value = 9257
$ws.Range("K132").Value = 8729.045999999998
$ws.Range("L132").Value = 27771
$ws.Range("M132").Value = -6199.045999999998
$ws.Range("N132").Value = -32831

$ws.Range("H140").Value = 69999.17999999999
$ws.Range("J140").Value = 69999.17999999999
$ws.Range("L140").Value = 69999.17999999999
$ws.Range("N140").Value = -80359.17999999999

$ws.Range("H141").Value = 159355.78
$ws.Range("J141").Value = 159355.78
$ws.Range("L141").Value = 159355.78
$ws.Range("N141").Value = -169715.78

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3637.5417
$ws.Range("I131").Value = 2815.9
$ws.Range("K131").Value = 8447.700000000001
$ws.Range("M131").Value = -3407.700000000001

$ws.Range("H137").Value = 1596.0714
$ws.Range("I137").Value = 1377
$ws.Range("J137").Value = 2399.3333
$ws.Range("K137").Value = 4131
$ws.Range("L137").Value = 7197.999899999999
$ws.Range("M137").Value = 969
$ws.Range("N137").Value = -17397.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 98326.164
$ws.Range("I70").Value = 162572
$ws.Range("K70").Value = 162572
$ws.Range("M70").Value = -162302

$ws.Range("H73").Value = 98326.164
$ws.Range("I73").Value = 162572
$ws.Range("K73").Value = 162572
$ws.Range("M73").Value = -161636

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3134.3462
$ws.Range("I7").Value = 1655.5
$ws.Range("J7").Value = 5500.5
$ws.Range("K7").Value = 1655.5
$ws.Range("L7").Value = 5500.5
$ws.Range("M7").Value = -1543.5
$ws.Range("N7").Value = -5724.5

$ws.Range("H40").Value = 7130
$ws.Range("I40").Value = 7370.3687
$ws.Range("J40").Value = 6844.5625
$ws.Range("K40").Value = 7370.3687
$ws.Range("L40").Value = 6844.5625
$ws.Range("M40").Value = -7234.3687
$ws.Range("N40").Value = -7116.5625

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H125").Value = 67857.30499999999
$ws.Range("J125").Value = 67857.30499999999
$ws.Range("L125").Value = 67857.30499999999
$ws.Range("N125").Value = -77697.30499999999

$ws.Range("H126").Value = 3134.3462
$ws.Range("I126").Value = 1655.5
$ws.Range("J126").Value = 5500.5
$ws.Range("K126").Value = 4966.5
$ws.Range("L126").Value = 16501.5
$ws.Range("M126").Value = -2496.5
$ws.Range("N126").Value = -21441.5

$ws.Range("H136").Value = 4048.9302
$ws.Range("I136").Value = 1959.8334
$ws.Range("K136").Value = 5879.5002
$ws.Range("M136").Value = -3329.5002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2346.2407
$ws.Range("I132").Value = 2132.1633
$ws.Range("K132").Value = 6396.4899
$ws.Range("M132").Value = -3866.4899

$ws.Range("H136").Value = 2463.9167
$ws.Range("I136").Value = 1220.4445
$ws.Range("K136").Value = 3661.3335
$ws.Range("M136").Value = -1111.3335
